$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 62
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = 13
$ws.Range("E42").Value = 20
$ws.Range("F42").Value = 76
$ws.Range("G42").Value = 96
